$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Simulation Time")

# --- New -O0 ("PH(0)") raw sample data -------------------------------------

# Row 2 ("TSN" row, formerly "TTS") - Simulink columns (I:R) have no data yet,
# but the Piha (O2) [T:AC] and Piha (O0) [AE:AN] columns get their first
# results.
$ws.Range("T2").Value = 126
$ws.Range("U2").Value = 124
$ws.Range("V2").Value = 128
$ws.Range("W2").Value = 115
$ws.Range("X2").Value = 115
$ws.Range("Y2").Value = 123
$ws.Range("Z2").Value = 115
$ws.Range("AA2").Value = 128
$ws.Range("AB2").Value = 102
$ws.Range("AC2").Value = 112

$ws.Range("AE2").Value = 313
$ws.Range("AF2").Value = 321
$ws.Range("AG2").Value = 333
$ws.Range("AH2").Value = 314
$ws.Range("AI2").Value = 319
$ws.Range("AJ2").Value = 311
$ws.Range("AK2").Value = 312
$ws.Range("AL2").Value = 310
$ws.Range("AM2").Value = 312
$ws.Range("AN2").Value = 316

# Row 5 ("MTG") - Piha (O0) [AE:AN] results.
$ws.Range("AE5").Value = 252
$ws.Range("AF5").Value = 252
$ws.Range("AG5").Value = 249
$ws.Range("AH5").Value = 252
$ws.Range("AI5").Value = 251
$ws.Range("AJ5").Value = 244
$ws.Range("AK5").Value = 252
$ws.Range("AL5").Value = 243
$ws.Range("AM5").Value = 251
$ws.Range("AN5").Value = 242

# Row 6 ("NP") - Piha (O0) [AE:AN] results.
$ws.Range("AE6").Value = 506
$ws.Range("AF6").Value = 516
$ws.Range("AG6").Value = 511
$ws.Range("AH6").Value = 501
$ws.Range("AI6").Value = 487
$ws.Range("AJ6").Value = 492
$ws.Range("AK6").Value = 502
$ws.Range("AL6").Value = 506
$ws.Range("AM6").Value = 515
$ws.Range("AN6").Value = 499

# --- Rename row 2's label from "TTS" to "TSN" ------------------------------
$ws.Range("A2").Value = "TSN"

# --- Format the summary columns (B:D) as whole numbers ---------------------
$ws.Range("B2:D6").NumberFormat = "0"

# --- Column D (Piha O0 average) needs a touch-up to its width --------------
$ws.Columns.Item(4).ColumnWidth = 9.140625

# --- Misc page / view bookkeeping ------------------------------------------
$ws.PageSetup.Orientation = 1
$ws.Range("C14").Select()
